$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text columns to be stored as text (not auto-converted to numbers/dates)
$ws.Range("A1:B6").NumberFormat = "@"
$ws.Range("E1:F6").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "Person"
$ws.Range("B1").Value = "Category"
$ws.Range("C1").Value = "Expense"
$ws.Range("D1").Value = "Savings"
$ws.Range("E1").Value = "Date"
$ws.Range("F1").Value = "Description"

# Row 2 - father / Investment
$ws.Range("A2").Value = "father"
$ws.Range("B2").Value = "Investment"
$ws.Range("C2").Value = 25000
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "2026-02-26"
$ws.Range("F2").Value = ""

# Row 3 - mom / Groceries
$ws.Range("A3").Value = "mom"
$ws.Range("B3").Value = "Groceries"
$ws.Range("C3").Value = 2500
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "2026-02-20"
$ws.Range("F3").Value = "food"

# Row 4 - son / Education
$ws.Range("A4").Value = "son "
$ws.Range("B4").Value = "Education"
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "2026-02-20"
$ws.Range("F4").Value = "bus ticket"

# Row 5 - father / Savings
$ws.Range("A5").Value = "father "
$ws.Range("B5").Value = "Savings"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1000
$ws.Range("E5").Value = "2026-02-21"
$ws.Range("F5").Value = ""

# Row 6 - TOTAL
$ws.Range("A6").Value = "TOTAL"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = 27600
$ws.Range("D6").Value = 1000
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
